$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("K2").Value = 2470
$ws.Range("K3").Value = 2376
$ws.Range("K4").Value = 500
$ws.Range("K6").Value = 2963
$ws.Range("K7").Value = 8463

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("K6").Value = 67
$ws.Range("K7").Value = 255
$ws.Range("K8").Value = 566
$ws.Range("K11").Value = 183
$ws.Range("K12").Value = 12
$ws.Range("K14").Value = 48
$ws.Range("K15").Value = 84
$ws.Range("K16").Value = 24
$ws.Range("K18").Value = 56
$ws.Range("K19").Value = 250
$ws.Range("K20").Value = 188
$ws.Range("K23").Value = 78
$ws.Range("K29").Value = 433
$ws.Range("K33").Value = 327
$ws.Range("K37").Value = 271
$ws.Range("K42").Value = 297
$ws.Range("K46").Value = 17
$ws.Range("K47").Value = 47
$ws.Range("K48").Value = 103
$ws.Range("K49").Value = 56
$ws.Range("K50").Value = 54
$ws.Range("K51").Value = 92
$ws.Range("K52").Value = 232
$ws.Range("K53").Value = 123
$ws.Range("K54").Value = 160
$ws.Range("K55").Value = 93
$ws.Range("K63").Value = 32
$ws.Range("K64").Value = 55
$ws.Range("K65").Value = 201
$ws.Range("K74").Value = 10
$ws.Range("K76").Value = 121
$ws.Range("K78").Value = 112
$ws.Range("K79").Value = 220
$ws.Range("K83").Value = 187
$ws.Range("K84").Value = 61
$ws.Range("K85").Value = 410
$ws.Range("K86").Value = 55
$ws.Range("K94").Value = 103
$ws.Range("K95").Value = 135
$ws.Range("K99").Value = 152
$ws.Range("K101").Value = 8463

$ws = $wb.Worksheets.Item("Bridgeport")
$ws.Range("K2").Value = 19
$ws.Range("K7").Value = 48

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("K2").Value = 84
$ws.Range("K7").Value = 255

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("K2").Value = 56
$ws.Range("K7").Value = 183

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("K2").Value = 150
$ws.Range("K3").Value = 139
$ws.Range("K7").Value = 410

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("K3").Value = 54
$ws.Range("K7").Value = 232

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("K6").Value = 63
$ws.Range("K7").Value = 123

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("K2").Value = 167
$ws.Range("K3").Value = 166
$ws.Range("K7").Value = 566

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("K4").Value = 12
$ws.Range("K7").Value = 187

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("K2").Value = 93
$ws.Range("K4").Value = 19
$ws.Range("K6").Value = 89
$ws.Range("K7").Value = 327

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("K3").Value = 43
$ws.Range("K6").Value = 38
$ws.Range("K7").Value = 135

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("K2").Value = 69
$ws.Range("K3").Value = 96
$ws.Range("K7").Value = 271

$ws = $wb.Worksheets.Item("New City")
$ws.Range("K2").Value = 59
$ws.Range("K3").Value = 51
$ws.Range("K7").Value = 201

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("K3").Value = 57
$ws.Range("K7").Value = 152

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("K6").Value = 17
$ws.Range("K7").Value = 61

$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("K4").Value = 6
$ws.Range("K6").Value = 35
$ws.Range("K7").Value = 56

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("K6").Value = 69
$ws.Range("K7").Value = 160

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("K2").Value = 115
$ws.Range("K6").Value = 141
$ws.Range("K7").Value = 433

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("K3").Value = 17
$ws.Range("K6").Value = 53
$ws.Range("K7").Value = 103

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("K3").Value = 66
$ws.Range("K6").Value = 87
$ws.Range("K7").Value = 250

$ws = $wb.Worksheets.Item("River North")
$ws.Range("K2").Value = 22
$ws.Range("K6").Value = 72
$ws.Range("K7").Value = 121

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("K2").Value = 23
$ws.Range("K6").Value = 20
$ws.Range("K7").Value = 67

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("K3").Value = 91
$ws.Range("K6").Value = 119
$ws.Range("K7").Value = 297

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("K6").Value = 42
$ws.Range("K7").Value = 112

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("K3").Value = 21
$ws.Range("K7").Value = 93

$ws = $wb.Worksheets.Item("Jefferson Park")
$ws.Range("K6").Value = 7
$ws.Range("K7").Value = 17

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("K6").Value = 23
$ws.Range("K7").Value = 78

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("K2").Value = 72
$ws.Range("K3").Value = 78
$ws.Range("K6").Value = 51
$ws.Range("K7").Value = 220

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("K6").Value = 18
$ws.Range("K7").Value = 55

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("K2").Value = 64
$ws.Range("K3").Value = 53
$ws.Range("K7").Value = 188

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("K6").Value = 12
$ws.Range("K7").Value = 56

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("K6").Value = 47
$ws.Range("K7").Value = 103

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("K3").Value = 17
$ws.Range("K7").Value = 47

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("K2").Value = 28
$ws.Range("K7").Value = 84

$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Range("K4").Value = 6
$ws.Range("K7").Value = 54

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("K3").Value = 11
$ws.Range("K7").Value = 55

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("K2").Value = 23
$ws.Range("K3").Value = 25
$ws.Range("K4").Value = 10
$ws.Range("K7").Value = 92

$ws = $wb.Worksheets.Item("Beverly")
$ws.Range("K6").Value = 3
$ws.Range("K7").Value = 12

$ws = $wb.Worksheets.Item("Bucktown")
$ws.Range("K6").Value = 15
$ws.Range("K7").Value = 24

$ws = $wb.Worksheets.Item("Printers Row")
$ws.Range("K6").Value = 7
$ws.Range("K7").Value = 10
